# Scheduled-runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit
# tables on each job sheet. Values only - no structural / formatting change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 32783.332
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 32783.332
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 32783.332
$ws.Range("N21").Value = -33719.332

$ws.Range("H23").Value = 32783.332
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 32783.332
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 32783.332
$ws.Range("N23").Value = -33251.332

$ws.Range("H34").Value = 952.8
$ws.Range("I34").Value = 952.8
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 952.8
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -749.8
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 952.8
$ws.Range("I36").Value = 952.8
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 952.8
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -237.8
$ws.Range("N36").ClearContents()

$ws.Range("H68").Value = 18750
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 18750
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 18750
$ws.Range("N68").Value = -20248

$ws.Range("H70").Value = 1542.48
$ws.Range("I70").Value = 1315.5
$ws.Range("J70").Value = 1585.7142
$ws.Range("K70").Value = 3946.5
$ws.Range("L70").Value = 4757.142599999999
$ws.Range("M70").Value = -3676.5
$ws.Range("N70").Value = -5297.142599999999

$ws.Range("H71").Value = 18750
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 18750
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 56250
$ws.Range("N71").Value = -63738

$ws.Range("H73").Value = 1542.48
$ws.Range("I73").Value = 1315.5
$ws.Range("J73").Value = 1585.7142
$ws.Range("K73").Value = 3946.5
$ws.Range("L73").Value = 4757.142599999999
$ws.Range("M73").Value = -3010.5
$ws.Range("N73").Value = -6629.142599999999

$ws.Range("H100").Value = 2046.6666
$ws.Range("I100").Value = 908
$ws.Range("J100").Value = 3185.3333
$ws.Range("K100").Value = 908
$ws.Range("L100").Value = 3185.3333
$ws.Range("M100").Value = -367
$ws.Range("N100").Value = -4267.3333

$ws.Range("H137").Value = 8000879
$ws.Range("I137").Value = 868.05554
$ws.Range("J137").Value = 28572336
$ws.Range("K137").Value = 2604.16662
$ws.Range("L137").Value = 85717008
$ws.Range("M137").Value = -54.16661999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8652.606
$ws.Range("I32").Value = 8217.053
$ws.Range("J32").Value = 11411.111
$ws.Range("K32").Value = 8217.053
$ws.Range("L32").Value = 11411.111
$ws.Range("M32").Value = -7930.053
$ws.Range("N32").Value = -11985.111

$ws.Range("H132").Value = 7577375
$ws.Range("I132").Value = 9616622
$ws.Range("J132").Value = 3028.2856
$ws.Range("K132").Value = 28849866
$ws.Range("L132").Value = 9084.856800000001
$ws.Range("M132").Value = -28847336
$ws.Range("N132").Value = -14144.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 43927
$ws.Range("I75").Value = 22500
$ws.Range("J75").Value = 48212.4
$ws.Range("K75").Value = 22500
$ws.Range("L75").Value = 48212.4
$ws.Range("M75").Value = -21564
$ws.Range("N75").Value = -50084.4

$ws.Range("H78").Value = 43927
$ws.Range("I78").Value = 22500
$ws.Range("J78").Value = 48212.4
$ws.Range("K78").Value = 67500
$ws.Range("L78").Value = 144637.2
$ws.Range("M78").Value = -62820
$ws.Range("N78").Value = -153997.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5851157.5
$ws.Range("I31").Value = 3234.2036
$ws.Range("J31").Value = 111113784
$ws.Range("K31").Value = 3234.2036
$ws.Range("L31").Value = 111113784
$ws.Range("M31").Value = -2939.2036
$ws.Range("N31").Value = -111114374

$ws.Range("H34").Value = 5851157.5
$ws.Range("I34").Value = 3234.2036
$ws.Range("J34").Value = 111113784
$ws.Range("K34").Value = 3234.2036
$ws.Range("L34").Value = 111113784
$ws.Range("M34").Value = -3032.2036
$ws.Range("N34").Value = -111114188

$ws.Range("H62").Value = 2542.7856
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 3159.8
$ws.Range("K62").Value = 2200
$ws.Range("L62").Value = 3159.8
$ws.Range("M62").Value = -1576
$ws.Range("N62").Value = -4407.8

$ws.Range("H65").Value = 2542.7856
$ws.Range("I65").Value = 2200
$ws.Range("J65").Value = 3159.8
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 15799
$ws.Range("M65").Value = -7880
$ws.Range("N65").Value = -22039

$ws.Range("H70").Value = 50090
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 50090
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 50090
$ws.Range("N70").Value = -50720

$ws.Range("H73").Value = 50090
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 50090
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 50090
$ws.Range("N73").Value = -52274

$ws.Range("H132").Value = 26318462
$ws.Range("I132").Value = 41668430
$ws.Range("J132").Value = 4230.2856
$ws.Range("K132").Value = 125005290
$ws.Range("L132").Value = 12690.8568
$ws.Range("M132").Value = -125002760
$ws.Range("N132").Value = -17750.8568

$ws.Range("H140").Value = 44459.832
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 44459.832
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 44459.832
$ws.Range("N140").Value = -54819.832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1225
$ws.Range("I117").Value = 199
$ws.Range("J117").Value = 1430.2
$ws.Range("K117").Value = 597
$ws.Range("L117").Value = 4290.6
$ws.Range("M117").Value = 2845
$ws.Range("N117").Value = -11174.6

$ws.Range("H131").Value = 825.9400000000001
$ws.Range("I131").Value = 388.2857
$ws.Range("J131").Value = 858.8817
$ws.Range("K131").Value = 1164.8571
$ws.Range("L131").Value = 2576.6451
$ws.Range("M131").Value = 3875.1429
$ws.Range("N131").Value = -12656.6451

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4072.5312
$ws.Range("I126").Value = 2841.4614
$ws.Range("J126").Value = 4914.8423
$ws.Range("K126").Value = 8524.3842
$ws.Range("L126").Value = 14744.5269
$ws.Range("M126").Value = -6054.3842

$ws.Range("H132").Value = 4180.269
$ws.Range("I132").Value = 2893.4707
$ws.Range("J132").Value = 6610.8887
$ws.Range("K132").Value = 8680.4121
$ws.Range("L132").Value = 19832.6661
$ws.Range("M132").Value = -6150.4121
$ws.Range("N132").Value = -24892.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7446.5835
$ws.Range("I40").Value = 8365.571
$ws.Range("J40").Value = 6160
$ws.Range("K40").Value = 8365.571
$ws.Range("L40").Value = 6160
$ws.Range("M40").Value = -8229.571
$ws.Range("N40").Value = -6432

$ws.Range("H61").Value = 2041.4286
$ws.Range("I61").Value = 1897.5
$ws.Range("J61").Value = 2233.3333
$ws.Range("K61").Value = 1897.5
$ws.Range("L61").Value = 2233.3333
$ws.Range("M61").Value = -1695.5

$ws.Range("H113").Value = 2041.4286
$ws.Range("I113").Value = 1897.5
$ws.Range("J113").Value = 2233.3333
$ws.Range("K113").Value = 1897.5
$ws.Range("L113").Value = 2233.3333
$ws.Range("M113").Value = 272.5

$ws.Range("H122").Value = 4885.25
$ws.Range("I122").Value = 4632.7617
$ws.Range("J122").Value = 5642.7144
$ws.Range("K122").Value = 13898.2851
$ws.Range("L122").Value = 16928.1432
$ws.Range("M122").Value = -11448.2851
$ws.Range("N122").Value = -21828.1432

$ws.Range("H132").Value = 13897658
$ws.Range("I132").Value = 7638
$ws.Range("J132").Value = 21748540
$ws.Range("K132").Value = 22914
$ws.Range("L132").Value = 65245620
$ws.Range("M132").Value = -20384
$ws.Range("N132").Value = -65250680

$ws.Range("H139").Value = 60663.57
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 60663.57
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 60663.57
$ws.Range("N139").Value = -70943.57000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2415.3257
$ws.Range("I96").Value = 1881.96
$ws.Range("J96").Value = 3156.111
$ws.Range("K96").Value = 1881.96
$ws.Range("L96").Value = 3156.111
$ws.Range("M96").Value = -508.96
$ws.Range("N96").Value = -5902.111

$ws.Range("H132").Value = 1427.7142
$ws.Range("I132").Value = 745.26666
$ws.Range("J132").Value = 3133.8333
$ws.Range("K132").Value = 2235.79998
$ws.Range("L132").Value = 9401.499899999999
$ws.Range("M132").Value = 294.2000200000002
$ws.Range("N132").Value = -14461.4999
